# Project Report.docx edit script
# 1) Paragraph 2 ("Patrick Bauer") -> split into 4 runs:
#      "David Babcock" | ", " | "Patrick Bauer" | ", Samuel Arminana, Maya Awad"
# 2) Paragraph 6 (database design) -> strip spell-check (proofErr) run-splits,
#      collapse into a single run with identical text.
# 3) Paragraph 9 (WAMPserver) -> strip spell-check run-splits for the first part
#      of the paragraph (collapse into one run) while leaving the trailing,
#      untouched "  To style our webpages..." run (rsid-tagged) intact/separate.
# 4) Paragraph 10 (mysqli_connect) -> strip spell-check run-splits, collapse
#      into a single run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: touch (identity find/replace) a bounded sub-range so the engine
# re-serialises it as a single run using the formatting of the first run in
# that span, merging away any <w:proofErr/> markers and run splits that
# carry identical formatting. Any live $d.Bookmarks act as merge barriers,
# so callers can bracket spans they want to keep separate.
# ---------------------------------------------------------------------------
function Touch-Range($rng) {
    $t = $rng.Text
    [void]$rng.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, $t, 2)
}

# ===========================================================================
# 1) Names paragraph
# ===========================================================================
$namesFind = $d.Content
[void]$namesFind.Find.Execute("Patrick Bauer", $true, $false, $false, $false, $false, `
    $true, 1, $false, "David Babcock, Patrick Bauer, Samuel Arminana, Maya Awad", 2)

$pNames = $d.Paragraphs(2)
$pStart = $pNames.Range.Start

# Desired run boundaries within "David Babcock, Patrick Bauer, Samuel Arminana, Maya Awad"
#   "David Babcock"                -> [0, 13)
#   ", "                           -> [13, 15)
#   "Patrick Bauer"                -> [15, 28)
#   ", Samuel Arminana, Maya Awad" -> [28, end)
$posAB = $pStart + 13
$posBC = $pStart + 15
$posCD = $pStart + 28

$d.Bookmarks.Add("zzAB", $d.Range($posAB, $posAB))
$d.Bookmarks.Add("zzBC", $d.Range($posBC, $posBC))
$d.Bookmarks.Add("zzCD", $d.Range($posCD, $posCD))

# One identity edit on the first segment forces the whole paragraph to
# re-split at the bookmark barriers while preserving run formatting.
$segEnd = $d.Bookmarks("zzAB").Range.End
Touch-Range $d.Range($pStart, $segEnd)

$d.Bookmarks("zzAB").Delete()
$d.Bookmarks("zzBC").Delete()
$d.Bookmarks("zzCD").Delete()

# ===========================================================================
# 2) Database design paragraph: collapse spell-check run splits
# ===========================================================================
$pDb = $d.Paragraphs(6)
Touch-Range $pDb.Range

# ===========================================================================
# 3) Software / WAMPserver paragraph: collapse spell-check run splits for the
#    first sentence span only; keep the trailing untouched run separate.
# ===========================================================================
$pSoft = $d.Paragraphs(9)
$tailFind = $d.Content
[void]$tailFind.Find.Execute("  To style our webpages", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$tailStart = $tailFind.Start
$d.Bookmarks.Add("zzTail", $d.Range($tailStart, $tailStart))

$softStart = $pSoft.Range.Start
$barrierPos = $d.Bookmarks("zzTail").Range.End
Touch-Range $d.Range($softStart, $barrierPos)

$d.Bookmarks("zzTail").Delete()

# ===========================================================================
# 4) mysqli_connect paragraph: collapse spell-check run splits
# ===========================================================================
$pConn = $d.Paragraphs(10)
Touch-Range $pConn.Range
